$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds the previous price-check snapshot (header timestamp in D1,
# then a price per product row, or blank for rows with no recorded price).
# This edit adds a new snapshot column E that starts out as a duplicate of
# column D (same values/types/formatting), then stamps the new run's
# timestamp into the E1 header.

$ws.Range("D1:D204").Copy($ws.Range("E1:E204"))
$ws.Range("E1").Value = "2026-01-27 18:21:51"
